$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 10 new columns starting at column AC (29), shifting existing
# columns AC..BK to AM..BU.
$insertRange = $ws.Range("AC1:AL1").EntireColumn
$insertRange.Insert()

# New header labels for the inserted columns (AC1:AL1)
$newHeaders = @("Negro","Blanco","Dorado","Plateado","Acero","Nude","Tonos marrones","Tonos pastel","Varios colores","Amarillo")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, 29 + $i).Value = $newHeaders[$i]
}

# Mark "Negro" (column AC) as available ("SI") for the product rows that
# already have color variants defined.
$ws.Range("AC2").Value = "SI"
$ws.Range("AC3").Value = "SI"
$ws.Range("AC4").Value = "SI"
$ws.Range("AC6").Value = "SI"

$ws.Range("AF9").Select()
